# Update odds values in rows 2, 3, 4, 5, 6 and 8 (these rows are above the
# row that gets removed, so updating them first or after the deletion makes
# no difference for their addressing).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("BD2").Value = 126

$ws.Range("I3").Value = 2.3

$ws.Range("G4").Value = 2.1

$ws.Range("G5").Value = 1.7

$ws.Range("G6").Value = 2.25
$ws.Range("N6").Value = 9.5
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.75
$ws.Range("X6").Value = 11
$ws.Range("AH6").Value = 15
$ws.Range("AT6").Value = 2.75
$ws.Range("AX6").Value = 17

$ws.Range("G8").Value = 2.35
$ws.Range("I8").Value = 2.92
$ws.Range("J8").Value = 2.92
$ws.Range("L8").Value = 3.55
$ws.Range("T8").Value = 2.55
$ws.Range("U8").Value = 1.6
$ws.Range("W8").Value = 8.75
$ws.Range("X8").Value = 13
$ws.Range("Z8").Value = 26
$ws.Range("AA8").Value = 18
$ws.Range("AB8").Value = 25
$ws.Range("AG8").Value = 9.25
$ws.Range("AH8").Value = 15.5
$ws.Range("AJ8").Value = 37
$ws.Range("AK8").Value = 25
$ws.Range("AL8").Value = 30
$ws.Range("AN8").Value = 4.3
$ws.Range("AP8").Value = 18.5
$ws.Range("AQ8").Value = 50
$ws.Range("AR8").Value = 75
$ws.Range("AS8").Value = 200
$ws.Range("AW8").Value = 4.9
$ws.Range("AX8").Value = 16.5
$ws.Range("AY8").Value = 23
$ws.Range("AZ8").Value = 80
$ws.Range("BA8").Value = 110
$ws.Range("BB8").Value = 300

# Remove row 9 (Charlotte vs Orlando City fixture) entirely; this shifts the
# two rows below (Colorado Rapids/LA Galaxy and Las Vegas Lights/Sacramento
# Republic) up by one, matching the new A1:BD10 dimension.
$ws.Rows.Item(9).Delete()
